# Update the cached "datetimeFigureOut" field text from 5/10/2011 to
# 5/11/2011 everywhere it appears: the Slide Master, every Slide Layout,
# and the Notes Master.

$p = $ppt.ActivePresentation

$oldDate = "5/10/2011"
$newDate = "5/11/2011"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }

        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
        } catch {
            $isDatePh = $false
        }
        if (-not $isDatePh) { continue }

        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout ("Custom Layout") hanging off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Notes Master (NB: $p.HasNotesMaster is unreliable in this host, so just
# go straight for $p.NotesMaster)
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
